$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '60.064.26'
$ws.Range('E2').Value = "'" + '  -0.76%  '
$ws.Range('D3').Value = "'" + '2.419.30'
$ws.Range('E3').Value = "'" + '  -1.25%  '
$ws.Range('E4').Value = "'" + '  +0.08%  '
$ws.Range('D5').Value = "'" + '552.63'
$ws.Range('E5').Value = "'" + '  -0.88%  '
$ws.Range('D6').Value = "'" + '137.24'
$ws.Range('E6').Value = "'" + '  -1.55%  '
$ws.Range('D7').Value = "'" + '1.00'
$ws.Range('E7').Value = "'" + '  +0.07%  '
$ws.Range('E8').Value = "'" + '  +3.69%  '
$ws.Range('D9').Value = "'" + '0.106'
$ws.Range('E9').Value = "'" + '  -1.67%  '
$ws.Range('D10').Value = "'" + '5.68'
$ws.Range('E10').Value = "'" + '  -2.43%  '
$ws.Range('D11').Value = "'" + '0.148'
$ws.Range('E11').Value = "'" + '  -0.95%  '
$ws.Range('E12').Value = "'" + '  -2.42%  '
$ws.Range('E13').Value = "'" + '  +0.57%  '
$ws.Range('D14').Value = "'" + '2.850.54'
$ws.Range('E14').Value = "'" + '  -0.99%  '
$ws.Range('D15').Value = "'" + '59.997.20'
$ws.Range('E16').Value = "'" + '  -2.03%  '
$ws.Range('D17').Value = "'" + '2.456.01'
$ws.Range('E17').Value = "'" + '  -0.34%  '
$ws.Range('D19').Value = "'" + '4.42'
$ws.Range('E19').Value = "'" + '  -0.62%  '
$ws.Range('D20').Value = "'" + '328.07'
$ws.Range('E20').Value = "'" + '  -2.42%  '
$ws.Range('D21').Value = "'" + '6.67'
$ws.Range('E21').Value = "'" + '  -3.61%  '
$ws.Range('E22').Value = "'" + '  +0.08%  '
$ws.Range('D23').Value = "'" + '65.93'
$ws.Range('E23').Value = "'" + '  +1.86%  '
$ws.Range('E24').Value = "'" + '  +2.76%  '
$ws.Range('D25').Value = "'" + '8.62'
$ws.Range('E25').Value = "'" + '  +0.56%  '
$ws.Range('E26').Value = "'" + '  -0.08%  '
$ws.Range('D27').Value = "'" + '1.40'
$ws.Range('E27').Value = "'" + '  +1.34%  '
$ws.Range('D28').Value = "'" + '0.0₃0777'
$ws.Range('E28').Value = "'" + '  -2.81%  '
$ws.Range('E29').Value = "'" + '  -2.43%  '
$ws.Range('D30').Value = "'" + '169.15'
$ws.Range('E30').Value = "'" + '  -1.06%  '
$ws.Range('D31').Value = "'" + '6.05'
$ws.Range('E31').Value = "'" + '  -4.32%  '
$ws.Range('E32').Value = "'" + '  +1.08%  '
$ws.Range('D33').Value = "'" + '18.58'
$ws.Range('E33').Value = "'" + '  -1.56%  '
$ws.Range('D35').Value = "'" + '1.32'
$ws.Range('E35').Value = "'" + '  -0.80%  '
$ws.Range('E36').Value = "'" + '  +0.04%  '
$ws.Range('E37').Value = "'" + '  -2.31%  '
$ws.Range('E38').Value = "'" + '  -2.24%  '
$ws.Range('D39').Value = "'" + '324.31'
$ws.Range('E39').Value = "'" + '  +2.07%  '
$ws.Range('E40').Value = "'" + '  -3.54%  '
$ws.Range('E41').Value = "'" + '  -2.25%  '
$ws.Range('D42').Value = "'" + '140.40'
$ws.Range('E42').Value = "'" + '  -2.87%  '
$ws.Range('D43').Value = "'" + '0.0971'
$ws.Range('E43').Value = "'" + '  +0.54%  '
$ws.Range('E44').Value = "'" + '  -1.97%  '
$ws.Range('D45').Value = "'" + '0.0516'
$ws.Range('E45').Value = "'" + '  -1.90%  '
$ws.Range('E46').Value = "'" + '  +0.37%  '
$ws.Range('D47').Value = "'" + '0.0223'
$ws.Range('E47').Value = "'" + '  -1.89%  '
$ws.Range('D48').Value = "'" + '0.387'
$ws.Range('E48').Value = "'" + '  -3.69%  '
$ws.Range('E49').Value = "'" + '  -0.06%  '
$ws.Range('E50').Value = "'" + '  -5.11%  '
$ws.Range('E51').Value = "'" + '  -1.03%  '
